# feat(transaction): transaction validation for role pegawai
# Insert a new "PTKP" column into the employee template (between NPWP and Agama).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Agama" and everything after shifts right by one)
$ws.Range("E1").EntireColumn.Insert()

# New header for the inserted column
$ws.Range("E1").Value = "PTKP"

# Re-apply column widths for the full A:O range to match the refreshed layout
$ws.Columns.Item(1).ColumnWidth = 13.166666666666666
$ws.Columns.Item(2).ColumnWidth = 16.166666666666668
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668
$ws.Columns.Item(4).ColumnWidth = 15.666666666666666
$ws.Columns.Item(5).ColumnWidth = 15.666666666666666
$ws.Columns.Item(6).ColumnWidth = 12.0
$ws.Columns.Item(7).ColumnWidth = 11.666666666666666
$ws.Columns.Item(8).ColumnWidth = 17.333333333333332
$ws.Columns.Item(9).ColumnWidth = 12.5
$ws.Columns.Item(10).ColumnWidth = 16.166666666666668
$ws.Columns.Item(11).ColumnWidth = 14.333333333333334
$ws.Columns.Item(12).ColumnWidth = 5.0
$ws.Columns.Item(13).ColumnWidth = 13.666666666666666
$ws.Columns.Item(14).ColumnWidth = 6.5
$ws.Columns.Item(15).ColumnWidth = 10.166666666666666

# Reset the view: scroll back to A1 and move the selection to G6
$ws.Range("A1").Select() | Out-Null
$ws.Range("G6").Select() | Out-Null
